$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New merged/centered header spanning G5:H5: "exchange time (s)" ---
$ws.Range("G5:H5").Merge()
$ws.Range("G5:H5").HorizontalAlignment = -4108  # xlCenter
$ws.Range("G5").Value = "exchange time (s)"

# --- Row 6 sub-headers for the two data series ---
$ws.Range("G6").Value = "23TAG"
$ws.Range("H6").Value = "T30177-TT"

# --- New column H values (second HDX-MS series) ---
$ws.Range("H7").Value = 3.3079570258276223
$ws.Range("H8").Value = 4.1349462822845275
$ws.Range("H9").Value = 4.7256528940394604
$ws.Range("H10").Value = 5.5132617097127037
$ws.Range("H11").Value = 11.073968936476858
$ws.Range("H12").Value = 12.655964498830695
$ws.Range("H13").Value = 14.765291915302477
$ws.Range("H14").Value = 25.718112194695514
$ws.Range("H15").Value = 29.392128222509157
$ws.Range("H16").Value = 34.290816259594017
$ws.Range("H17").Value = 53.821039482092587
$ws.Range("H18").Value = 61.509759408105815
$ws.Range("H19").Value = 71.761385976123449
$ws.Range("H20").Value = 80.315545955149929
$ws.Range("H21").Value = 91.789195377314215
$ws.Range("H22").Value = 107.08739460686658
$ws.Range("H23").Value = 136.13097058327409
$ws.Range("H24").Value = 163.35716469992892
$ws.Range("H25").Value = 201.49860328019528
$ws.Range("H26").Value = 235.08170382689448

# Two-decimal display for the new numeric column
$ws.Range("H7:H26").NumberFormat = "0.00"

# Column H width (close to the workbook's saved best-fit width)
$ws.Columns("H").ColumnWidth = 8.6666666666667

# Page set up to portrait (matches saved print settings)
$ws.PageSetup.Orientation = 1

# Final selection, as left by the editing session
$ws.Range("K12").Select()
